$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each coin row.
# D values are assigned with a leading apostrophe + Style reset so that
# numeric-looking strings (e.g. "0.4300", "80.50") are kept as literal
# text (preserving trailing zeros) instead of being coerced to numbers,
# while still ending up with the default (un-styled) cell format.
$ws.Range("D2").Value = "'27.506.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "'1.831.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'313.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "'0.4300"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.25%  "
$ws.Range("D8").Value = "'0.3664"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.27%  "
$ws.Range("D9").Value = "'0.07269"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.00%  "
$ws.Range("D10").Value = "'0.8677"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.14%  "
$ws.Range("D11").Value = "'20.64"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.59%  "
$ws.Range("D12").Value = "'1.865.67"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.35%  "
$ws.Range("D13").Value = "'5.407"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.25%  "
$ws.Range("D14").Value = "'6.533"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").Value = "'0.06936"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").Value = "'1.004"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").Value = "'80.50"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").Value = "'0.000008900"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.14%  "
$ws.Range("D20").Value = "'15.40"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("D21").Value = "'27.441.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.58%  "
$ws.Range("D22").Value = "'5.137"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.22%  "
$ws.Range("D23").Value = "'10.81"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.09%  "
$ws.Range("D24").Value = "'2.111.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.04%  "
$ws.Range("D25").Value = "'1.978"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "'154.44"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.08%  "
$ws.Range("D27").Value = "'18.82"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.00%  "
$ws.Range("D28").Value = "'5.147"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.97%  "
$ws.Range("D29").Value = "'114.06"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.04%  "
$ws.Range("D30").Value = "'1.829"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.71%  "
$ws.Range("D31").Value = "'0.08891"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("D32").Value = "'0.7536"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("D33").Value = "'2.986"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.48%  "
$ws.Range("D34").Value = "'4.543"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.04%  "
$ws.Range("D35").Value = "'1.135"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.99%  "
$ws.Range("D36").Value = "'1.001"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").Value = "'1.092"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.31%  "
$ws.Range("D38").Value = "'0.05323"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.70%  "
$ws.Range("D39").Value = "'0.01936"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.31%  "
$ws.Range("D40").Value = "'2.800"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.29%  "
$ws.Range("D41").Value = "'0.1667"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.57%  "
$ws.Range("D42").Value = "'0.5073"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.20%  "
$ws.Range("D43").Value = "'6.601"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.99%  "
$ws.Range("D44").Value = "'8.381"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.45%  "
$ws.Range("D45").Value = "'10.48"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.68%  "
$ws.Range("D46").Value = "'105.89"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.38%  "
$ws.Range("D47").Value = "'0.06491"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("D48").Value = "'0.4686"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("D49").Value = "'1.001"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("D50").Value = "'1.610"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.74%  "
$ws.Range("D51").Value = "'64.16"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.18%  "
